$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set new cell text values in the same order the original edit introduced the
# new shared strings, so the rebuilt shared-string table lines up with the
# target ordering as closely as possible.
$ws.Range("D1").Value = 'id="get_motorcycle"'
$ws.Range("C1").Value = 'id="get_camper"'
$ws.Range("F1").Value = 'id="get_automobile"'
$ws.Range("E1").Value = 'id="get_truck"'
$ws.Range("G3").Value = "npaInsurances"
$ws.Range("G1").Value = '//ol[@class="flex-control-nav flex-control-paging"]'
$ws.Range("H1").Value = 'id="read_more"'
$ws.Range("I1").Value = 'id="offer_automobile"'
$ws.Range("J1").Value = 'id="offer_camper"'
$ws.Range("K1").Value = 'id="offer_truck"'
$ws.Range("L1").Value = 'id="offer_motorcycle"'
$ws.Range("I4").Value = "<ENABLED>"
$ws.Range("J4").Value = "<ENABLED>"
$ws.Range("K4").Value = "<ENABLED>"
$ws.Range("L4").Value = "<ENABLED>"

# --- Row 2 fill (new yellow highlight style across C2:L2) ---
$ws.Range("C2:L2").Interior.Color = 65535

# --- Row 3 button/control names (re-ordered & renamed from the old rot* labels) ---
$ws.Range("C3").Value = "butCamperGetQuote"
$ws.Range("D3").Value = "butMotorcycleGetQuote"
$ws.Range("E3").Value = "butTruckGetQuote"
$ws.Range("F3").Value = "butAutomobileGetQuote"
$ws.Range("H3").Value = "lnkReadMore"
$ws.Range("I3").Value = "picAutomobile"
$ws.Range("J3").Value = "picCamper"
$ws.Range("K3").Value = "picTruck"
$ws.Range("L3").Value = "picMotorcycle"

# --- Drop the old M:P columns (their content either moves or goes away) so
#     they lose their stale custom widths, then re-create the single "Action"
#     tail column fresh in M. ---
$ws.Range("M1:P1").EntireColumn.Delete()
$ws.Range("M1").Value = "Action"
$ws.Range("M2").Value = "Action"
$ws.Range("M3").Value = "Action"
$ws.Range("M4").Value = "<NOP>"

# --- Column widths (C:L) to match the widened / re-proportioned table ---
$ws.Columns.Item(3).ColumnWidth = 22.0
$ws.Columns.Item(4).ColumnWidth = 25.833333333333332
$ws.Columns.Item(5).ColumnWidth = 21.5
$ws.Columns.Item(6).ColumnWidth = 23.166666666666668
$ws.Columns.Item(7).ColumnWidth = 36.666666666666664
$ws.Columns.Item(8).ColumnWidth = 19.666666666666668
$ws.Columns.Item(9).ColumnWidth = 19.666666666666668
$ws.Columns.Item(10).ColumnWidth = 19.666666666666668
$ws.Columns.Item(11).ColumnWidth = 19.666666666666668
$ws.Columns.Item(12).ColumnWidth = 19.666666666666668

# --- Selection state on save ---
$ws.Range("A4:XFD4").Select()
